$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the RANDBETWEEN()/shared-formula column (B2:B12) with its
# last-computed static values (formulas -> plain numbers).
$values = @{
    2  = 2731
    3  = 2262
    4  = 2479
    5  = 2144
    6  = 2845
    7  = 2403
    8  = 2584
    9  = 2147
    10 = 2224
    11 = 2112
    12 = 2149
}

foreach ($row in 2..12) {
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = $values[$row]
    # Mirror the explicit (blank/default) style the cells picked up when
    # Excel converted them from formulas to literal values.
    $cell.NumberFormat = "General"
}

# Move the sheet's selection from E10 to E11, as recorded in the saved view.
$ws.Range("E11").Select() | Out-Null
